# Insert 10 new rows above the existing data, then populate them with
# labels/units/types paired data, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 8 rows (header + 7 data rows) down by 10 rows so they
# land in rows 11-18, by inserting 10 new blank rows above row 1.
$ws.Range("A1:D10").EntireRow.Insert()

# New paired-data rows 1-10 (columns A and B only)
$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "a"

$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = "b"

$ws.Range("A3").Value = "C"
$ws.Range("B3").Value = "c"

$ws.Range("A4").Value = "D"

$ws.Range("A5").Value = "E"
$ws.Range("B5").Value = "e"

$ws.Range("A6").Value = "F"
$ws.Range("B6").Value = "f"

$ws.Range("A7").Value = "Unit 1"
$ws.Range("B7").Value = "u1"

$ws.Range("A8").Value = "Unit 2"
$ws.Range("B8").Value = "u2"

$ws.Range("A9").Value = "Type 1"
$ws.Range("B9").Value = "t1"

$ws.Range("A10").Value = "Type 2"
$ws.Range("B10").Value = "t2"

# Update the active selection / window view to match the target state.
$ws.Range("F17").Select()

$excel.ActiveWindow.Left = 75
$excel.ActiveWindow.Top = 345
